# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header/data layout of the sheet (columns A-H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 (bold, centered, bordered style)
# onto the new header cells I1 and J1 before/after setting their values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), row 2 through row 10.
$dataI = @(5, 8, 7, 5, 2, 6, 5, 6, 7)
$dataJ = @(5, 8, 7, 5, 2, 7, 5, 6, 7)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}

Write-Output "Added columns I0 and J0 (I and J) to the sheet"
